# dataset.xlsx: "Completed annotations, added comparison for minmax normalization"
#
# The author widened the Environmental Perfomance Index (col E) and Consumer
# Price Index (col G) columns so the newly-completed annotation/comparison
# text fits, and left the window scrolled near the bottom of the table with
# cell E107 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns E and G to fit the completed annotations / minmax-normalization
# comparison text (width is expressed in characters, like the Excel UI's
# Format > Column Width dialog).
$ws.Columns.Item(5).ColumnWidth = 22.5
$ws.Columns.Item(7).ColumnWidth = 27.333333333333332

# Leave the view scrolled down near the bottom of the data with E107 selected,
# matching where the author ended up after finishing the annotations.
$ws.Application.Goto($ws.Range("A77")) | Out-Null
$ws.Range("E107").Select() | Out-Null
